# Error Calculations and Plots
# Delete the "RM 232" row and the "SC 92" row from the data sheet.
# In the original sheet these are row 26 (RM 232) and row 28 (SC 92).
# Deleting row 26 first shifts SC 92 up to row 27, so it is deleted next.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()
